$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: update extraction timestamp
$ws.Range("C2").Value = "2025-01-23 16:22:41"

# Row 3: update extraction timestamp
$ws.Range("C3").Value = "2025-01-23 16:22:55"

# Row 4: update extraction timestamp and fill in municipio/estado
$ws.Range("C4").Value = "2025-01-23 16:23:01"
$ws.Range("D4").Value = "Cordeiro"
$ws.Range("E4").Value = "Rio de Janeiro"

# Row 93: update extraction timestamp
$ws.Range("C93").Value = "2025-01-23 16:23:08"

# Row 94: status changed to ERRO, update extraction timestamp
$ws.Range("B94").Value = "ERRO"
$ws.Range("C94").Value = "2025-01-23 16:23:14"
